# Add schematics pinout tables for V1.1 and V1.2 connectors to all three
# sheets (Tabelle1, Tabelle2, Tabelle3). Each sheet gets a new 26-row,
# 2-column (pin number / signal name) table appended below the existing
# pinout table.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Tabelle1 (sheet1) - new rows 32..57
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws1Signals = @(
    "KL30", "KL15", "KL58b", "TKML", "Reverse", "GND",
    "MFA_Hebel", "MFA_Hebel", "MFA_Hebel", "MFA_Hebel", "GND",
    "Innetemp", "GetrTemp", "StartstopBat", "Solar", "Zweitbat",
    "Starterbat", "OilTemp", "Außentemp", "GND",
    "SDA", "SCL", "GND", "CANH", "CANL", "GND"
)
for ($i = 0; $i -lt $ws1Signals.Length; $i++) {
    $row = 32 + $i
    $ws1.Cells.Item($row, 1).Value = $i + 1
    $ws1.Cells.Item($row, 2).Value = $ws1Signals[$i]
}

# ---------------------------------------------------------------------
# Tabelle2 (sheet2) - new rows 33..58
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Tabelle2")
$ws2Signals = @(
    "KL30", "KL15", "KL58b", "TKML", "K-Line", "GND",
    "MFA_Hebel", "MFA_Hebel", "MFA_Hebel", "MFA_Hebel", "GND",
    "aux_temp", "oiltemp", "Gearbox_temp", "Manifold", "Starterbat",
    "Zweitbat", "Solar", "D+/nc", "GND",
    "SDA", "SCL", "GND", "CANH", "CANL", "GND"
)
for ($i = 0; $i -lt $ws2Signals.Length; $i++) {
    $row = 33 + $i
    $ws2.Cells.Item($row, 1).Value = $i + 1
    $ws2.Cells.Item($row, 2).Value = $ws2Signals[$i]
}

# ---------------------------------------------------------------------
# Tabelle3 (sheet3) - new rows 32..57
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Tabelle3")
$ws3Signals = @(
    "KL30", "KL15", "KL58b", "TKML", "K-Line", "GND",
    "MFA_Hebel", "MFA_Hebel", "MFA_Hebel", "MFA_Hebel", "GND",
    "aux_temp", "oiltemp", "Gearbox_temp", "Manifold", "Starterbat",
    "Zweitbat", "Solar", "StartstopBat", "GND",
    "SDA", "SCL", "GND", "CANH", "CANL", "GND"
)
for ($i = 0; $i -lt $ws3Signals.Length; $i++) {
    $row = 32 + $i
    $ws3.Cells.Item($row, 1).Value = $i + 1
    $ws3.Cells.Item($row, 2).Value = $ws3Signals[$i]
}

# ---------------------------------------------------------------------
# Restore / update the view state (scroll position + selection) on each
# sheet to reflect where the user ended up editing, finishing back on
# Tabelle1 so it remains the selected tab.
# ---------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("A33:B58").Select()

$ws3.Activate()
$ws3.Range("C:C").Select()

$ws1.Activate()
$ws1.Range("C57").Select()
